# Append the new daily mod-count row (row 55) to the ModCounts sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A55 looks like a date ("2026/01/04") but must stay literal text, matching
# the other Date-column cells (which are all stored as text, not real
# dates). Entering it with a leading apostrophe forces Excel to keep it as
# text instead of auto-converting it to a date serial number; resetting the
# cell's style back to "Normal" afterwards clears the transient
# quote-prefix marker that the apostrophe entry leaves behind.
$ws.Range("A55").Value = "'2026/01/04"
$ws.Range("A55").Style = "Normal"

$ws.Range("B55").Value = "逃离鸭科夫"
$ws.Range("C55").Value = 1137

# Match the centered alignment used by every other data row (A3:C54).
$ws.Range("A55:C55").HorizontalAlignment = -4108
$ws.Range("A55:C55").VerticalAlignment = -4108
